# Append 1 row at the end of the data range (row 9), replicating a new
# form submission similar to the rows above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 9

# Column C ("الكمية") holds a numeric-looking value ("2323") that must stay
# text (like the rest of the sheet, flagged via ignoredError
# numberStoredAsText), so force a text number format before writing it.
$ws.Cells.Item($row, 3).NumberFormat = "@"

# Column A is blank in the source row; use the text-prefix apostrophe so
# it is written out as an empty text cell rather than a truly empty one.
$ws.Cells.Item($row, 1).Value = "'"
$ws.Cells.Item($row, 2).Value = "أحمد شريم"
$ws.Cells.Item($row, 3).Value = "2323"
$ws.Cells.Item($row, 4).Value = "ايتا"
$ws.Cells.Item($row, 5).Value = "الرحلة 2"
$ws.Cells.Item($row, 6).Value = "C2"
$ws.Cells.Item($row, 7).Value = "NRC"
$ws.Cells.Item($row, 8).Value = "٠١‏/٠٥‏/٢٠٢٥ ٠٥:٢٢:٢٩ م"

# Drop back to the default (unstyled) cell format for the two cells we
# explicitly formatted, so the new row matches the rest of the sheet,
# which carries no explicit style index.
$ws.Cells.Item($row, 1).Style = "Normal"
$ws.Cells.Item($row, 3).Style = "Normal"
